$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-06-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-12 Thursday", 2) | Out-Null

# Update table cells (addressed by row/col to avoid ambiguity from duplicate values)
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "42÷9=4, 6"
$t.Cell(1,2).Range.Text = "68÷3=22, 2"
$t.Cell(1,3).Range.Text = "92÷7=13, 1"
$t.Cell(1,4).Range.Text = "16÷4=4, 0"
$t.Cell(1,5).Range.Text = "20÷7=2, 6"
$t.Cell(5,1).Range.Text = "77÷4=19, 1"
$t.Cell(5,2).Range.Text = "47÷5=9, 2"
$t.Cell(5,3).Range.Text = "44÷6=7, 2"
$t.Cell(5,4).Range.Text = "84÷8=10, 4"
$t.Cell(5,5).Range.Text = "99÷8=12, 3"
$t.Cell(9,1).Range.Text = "55÷9=6, 1"
$t.Cell(9,2).Range.Text = "29÷6=4, 5"
$t.Cell(9,3).Range.Text = "93÷7=13, 2"
$t.Cell(9,4).Range.Text = "68÷3=22, 2"
$t.Cell(9,5).Range.Text = "92÷3=30, 2"
$t.Cell(13,1).Range.Text = "56÷2=28, 0"
$t.Cell(13,2).Range.Text = "95÷9=10, 5"
$t.Cell(13,3).Range.Text = "74÷5=14, 4"
$t.Cell(13,4).Range.Text = "16÷7=2, 2"
$t.Cell(13,5).Range.Text = "77÷2=38, 1"
$t.Cell(17,1).Range.Text = "49÷5=9, 4"
$t.Cell(17,2).Range.Text = "66÷2=33, 0"
$t.Cell(17,3).Range.Text = "76÷4=19, 0"
$t.Cell(17,4).Range.Text = "24÷5=4, 4"
$t.Cell(17,5).Range.Text = "18÷8=2, 2"
